# Update countries & provincias Spain
#
# 1) Reorder three pairs/groups of countries in the "Pais" sheet so that the
#    text shown for a handful of existing rows changes (the row's underlying
#    rank/position stays put, but the country name displayed there swaps -
#    this mirrors a reordering of the shared-string table in the source
#    workbook). 2) Refresh the "Datos actualizados" timestamp. 3) Update the
#    numeric COVID figures (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for the touched rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 15:50"

# --- Country name swaps (rows keep their rank, the label changes) -----
# Grecia / India / Islandia / Panama -> Grecia / Islandia / India / Panama
$ws.Range("A44").Value = "Islandia"
$ws.Range("A45").Value = "India"

# Islas Feroe / Malta / Senegal / Ghana / Costa de Marfil
#   -> Islas Feroe / Ghana / Malta / Senegal / Costa de Marfil
$ws.Range("A95").Value = "Ghana"
$ws.Range("A96").Value = "Malta"
$ws.Range("A97").Value = "Senegal"

# Georgia / Kirguistan / Montenegro / Bolivia
#   -> Georgia / Montenegro / Kirguistan / Bolivia
$ws.Range("A115").Value = "Montenegro"
$ws.Range("A116").Value = "Kirguistan"

# --- Numeric data refresh ----------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 123828
$ws.Range("C4").Value = 250
$ws.Range("E4").Value = 118361

# Row 20: Noruega
$ws.Range("B20").Value = 4235
$ws.Range("C20").Value = 220
$ws.Range("E20").Value = 4203

# Row 43: Grecia
$ws.Range("E43").Value = 972
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 37

# Row 44: now Islandia
$ws.Range("B44").Value = 1020
$ws.Range("C44").Value = 57
$ws.Range("D44").Value = 114
$ws.Range("E44").Value = 904
$ws.Range("F44").Value = 19
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 2

# Row 45: now India
$ws.Range("B45").Value = 987
$ws.Range("D45").Value = 87
$ws.Range("E45").Value = 875
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 25

# Row 95: now Ghana
$ws.Range("B95").Value = 152
$ws.Range("C95").Value = 11
$ws.Range("E95").Value = 145
$ws.Range("H95").Value = 5

# Row 96: now Malta
$ws.Range("B96").Value = 151
$ws.Range("C96").Value = 2
$ws.Range("D96").Value = 2
$ws.Range("E96").Value = 149
$ws.Range("F96").Value = 1

# Row 97: now Senegal
$ws.Range("B97").Value = 142
$ws.Range("C97").Value = 12
$ws.Range("D97").Value = 27
$ws.Range("E97").Value = 115
$ws.Range("F97").Value = 0
$ws.Range("H97").Value = 0

# Row 115: now Montenegro
$ws.Range("B115").Value = 85
$ws.Range("C115").Value = 1
$ws.Range("F115").Value = 1
$ws.Range("H115").Value = 1

# Row 116: now Kirguistan
$ws.Range("C116").Value = 26
$ws.Range("E116").Value = 83
$ws.Range("F116").Value = 0
$ws.Range("H116").Value = 0
